# Adds spacing to section headings and job-title paragraphs to increase
# readability, and removes a redundant empty spacer paragraph before the
# PUBLICATIONS heading while giving that heading its own spacing.

$d = $word.ActiveDocument

# --- Section heading paragraphs (bottom-border "divider" headings) ---
# OBJECTIVE, TECHNICAL SKILLS, WORK EXPERIENCE: add SpaceAfter = 6pt (120 twips)
foreach ($headingText in @("OBJECTIVE", "TECHNICAL SKILLS", "WORK EXPERIENCE")) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.Trim()
        if ($t -eq $headingText) {
            $p.Format.SpaceAfter = 6
            break
        }
    }
}

# --- Job title paragraphs that start a new role: add SpaceBefore = 4pt (80 twips) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if (($t -match "^Data Scientist - 1, Krutrim") -or ($t -match "^Data Scientist - 1, Ola Cabs")) {
        $p.Format.SpaceBefore = 4
    }
}

# --- Remove the empty spacer paragraph right before "PUBLICATIONS" and give
#     the PUBLICATIONS heading its own before/after spacing ---
# NOTE: set the spacing *before* deleting the preceding paragraph - once the
# preceding paragraph's range is deleted, a Paragraph object obtained earlier
# by index no longer tracks the same paragraph, so do the formatting first.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()
    if ($t -eq "PUBLICATIONS & RESEARCH") {
        $p.Format.SpaceBefore = 8
        $p.Format.SpaceAfter = 6
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text.Trim() -eq "") {
            $prev.Range.Delete()
        }
        break
    }
}
